$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").Value = "2025/12/05 22:00"
$ws.Range("B76").Value = "12,025位本"
$ws.Range("C76").Value = "37位 広告・宣伝 (本)"
$ws.Range("D76").Value = "49位商業デザイン"
$ws.Range("E76").Value = "774位ビジネス実用本"
$ws.Range("F76").Value = "-"
$ws.Range("G76").Value = "-"
